# Applies crypto price/volume updates per commit "Updated cryptos list" diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.683.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.15%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.874.10"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.72%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.23%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'331.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +2.60%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.24%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4735"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +4.37%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3951"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +2.11%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'47.87"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.13%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.08053"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +1.77%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'1.028"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.53%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'22.07"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +2.96%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.879.77"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.93%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'5.965"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'7.151"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.23%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  +0.36%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.00001048"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.23%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'86.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.30%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.06650"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +2.14%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'17.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.74%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  +0.40%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'27.676.96"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +1.13%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'5.524"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.22%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'11.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.09%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.311"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +1.56%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'2.110.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +2.13%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'159.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +3.92%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +2.47%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'2.103"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +1.76%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'5.592"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +2.18%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'122.24"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +1.28%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.9747"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +4.28%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.09550"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +2.54%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.456"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -2.06%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'3.591"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.65%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'5.351"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.56%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.06105"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +1.94%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.02255"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.81%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'1.232"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.86%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'8.203"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.67%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.6032"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +1.98%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1903"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.79%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'10.27"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +1.45%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'1.268"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.56%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.5719"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +1.57%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'12.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +2.19%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'1.949"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.26%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'3.389"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.45%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'Quant"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'115.26"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +6.24%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "'Cronos"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'0.06868"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.29%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +14.35%  "
$ws.Range("E51").Style = "Normal"
